# Daily attendance processing - reorder "Recorded By" (column G) entries.
# Rule observed in the target dataset: within each comma-separated list of
# recorder names/emails, the literal entry "System" (exact case) is moved to
# the end of the list while the remaining entries keep their relative order.
# If a cell's list does not contain an exact "System" entry, the entries are
# sorted alphabetically instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$changedCount = 0

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if ($val -eq "") { continue }

    $parts = @($val -split ", ")
    if ($parts.Count -le 1) { continue }

    $others = @()
    $hasSystem = $false
    foreach ($p in $parts) {
        if ($p.Equals("System")) {
            $hasSystem = $true
        } else {
            $others += $p
        }
    }

    if ($hasSystem) {
        $others += "System"
        $newVal = $others -join ", "
    } else {
        $sorted = $parts | Sort-Object
        $newVal = $sorted -join ", "
    }

    if (-not $newVal.Equals($val)) {
        $cell.Value2 = $newVal
        $changedCount++
    }
}

Write-Output "Changed $changedCount cells in column G"
